$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the whole "Meta description: ..." paragraph that currently
#    follows the H1 title.
# ---------------------------------------------------------------------
$metaRange = $d.Content
$metaFound = $metaRange.Find.Execute("Meta description", $true, $false,
    $false, $false, $false, $true, 1, $false, "", 0)
if ($metaFound) {
    $metaRange.Expand(4) | Out-Null   # wdParagraph -> whole paragraph incl. mark
    $metaRange.Delete()
}

# ---------------------------------------------------------------------
# 2. Insert a new bold paragraph "Play Age of the Gods Medusa & Monsters
#    Free | Review" right before the closing "Prompt: ..." paragraph, and
#    replace that paragraph's own text with the review meta description
#    (keeping its italic formatting).
# ---------------------------------------------------------------------
$promptRange = $d.Content
$promptFound = $promptRange.Find.Execute("Prompt: Create a Cartoon Style",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($promptFound) {
    $promptRange.Expand(4) | Out-Null   # whole "Prompt: ..." paragraph
    $promptStart = $promptRange.Start
    $promptEnd = $promptRange.End

    # Insert a paragraph break right at the start of the "Prompt: ..."
    # paragraph - this pushes its content into a new paragraph that
    # follows, leaving a fresh empty paragraph in its place.
    $breakPoint = $d.Range($promptStart, $promptStart)
    $breakPoint.InsertParagraphBefore()

    # The newly created (still empty) paragraph spans [$promptStart,
    # $promptStart + 1) (just its paragraph mark); fill it in with the
    # bold heading text via a small well-formed OOXML fragment so no
    # stray formatting (e.g. the italics of the neighbouring "Prompt"
    # run) leaks in.
    $newParaRange = $d.Range($promptStart, $promptStart + 1)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Age of the Gods Medusa &amp; Monsters Free | Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $newParaRange.InsertXML($xml) | Out-Null

    # The "Prompt: ..." paragraph now starts right after the newly
    # inserted paragraph (shifted forward by the length of the inserted
    # text). Re-locate it with Find once more and replace its text,
    # leaving its paragraph mark / italic run formatting untouched.
    $promptRange2 = $d.Content
    $promptRange2.Find.Execute("Prompt: Create a Cartoon Style", $true,
        $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $promptRange2.Expand(4) | Out-Null
    $textOnly = $d.Range($promptRange2.Start, $promptRange2.End - 1)
    $textOnly.Text = "Read our review of Age of the Gods Medusa & Monsters online slot game and play for free. Learn about the special features, RTP, and winning potential."
}

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
